# Updated cryptos list on Wed May  1 04:56:59 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.341.04"
$ws.Range("E2").Value = "  -4.95%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.010.32"
$ws.Range("E3").Value = "  -5.11%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.89%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.04%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.002.26"
$ws.Range("E8").Value = "  -5.24%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -2.39%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.27%  "

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.92%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.63%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -7.36%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.98%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +0.48%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "3.506.63"
$ws.Range("E16").Value = "  -5.07%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.007.39"
$ws.Range("E17").Value = "  -5.22%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "60.322.26"
$ws.Range("E18").Value = "  -4.99%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -0.41%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.67%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -5.23%  "

# Row 22 - Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.675"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.20%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.63%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("E24").Value = "  -1.79%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.29%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.22%  "

# Row 27 - FirstDigitalUSD
$ws.Range("E27").Value = "  +0.08%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  -5.16%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  -4.36%  "

# Row 30 - RenderToken
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.18%  "

# Row 31 - NEARProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.37%  "

# Row 32 - EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.05%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0951"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.71%  "

# Row 34 - Mantle->Filecoin
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.43%  "

# Row 35 - Filecoin->Mantle
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.938"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.10%  "

# Row 36 - OKB
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.29%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  -15.02%  "

# Row 38 - Cosmos
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.42%  "

# Row 39 - PEPE
$ws.Range("D39").Value = "0.0₃0667"
$ws.Range("E39").Value = "  -8.94%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -8.42%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  -3.72%  "

# Row 42 - Bittensor
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "373.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.74%  "

# Row 43 - Maker
$ws.Range("D43").Value = "2.683.50"
$ws.Range("E43").Value = "  -3.82%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  -6.95%  "

# Row 45 - USDe
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46 - Monero
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "122.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.73%  "

# Row 47 - TheGraph
$ws.Range("E47").Value = "  -6.45%  "

# Row 48 - Fetch.AI
$ws.Range("E48").Value = "  -5.31%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -3.27%  "

# Row 50 - InjectiveProtocol
$ws.Range("E50").Value = "  -6.81%  "

# Row 51 - ThetaToken->Cronos
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.132"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.21%  "

